# Apply hybrid bold + color (2C3E50) highlighting to quantitative
# impact metrics (percentages, dollar amounts, large numbers) across
# the achievements / work-experience bullet paragraphs.
#
# Word's Font.Color takes a BGR-packed integer (0x00BBGGRR), so the
# hex "2C3E50" (R=2C, G=3E, B=50) must be submitted as 0x00503E2C.
function Bold-Metric($range, $searchText) {
    $fr = $range.Duplicate
    $ok = $fr.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($ok) {
        $fr.Font.Bold = $true
        $fr.Font.Color = 5258796
    }
    return $ok
}

$d = $word.ActiveDocument

# "Partner - Siege Analytics" bullets
# • Discovered systematic race coding errors ... from 23% to 64%
$p = $d.Paragraphs.Item(10)
Bold-Metric $p.Range "23%"
Bold-Metric $p.Range "64%"

# • Utilized advanced sampling methods ... from ±4.2% to ±2.1%, increasing
#   voter turnout prediction accuracy from 71% to 87% ...
$p = $d.Paragraphs.Item(12)
Bold-Metric $p.Range "±4.2%"
Bold-Metric $p.Range "±2.1%"
Bold-Metric $p.Range "71%"
Bold-Metric $p.Range "87%"

# • Trigonometric algorithm ... reduced mapping costs by 73.5%, saving
#   campaigns and organizations $4.7M ...
$p = $d.Paragraphs.Item(13)
Bold-Metric $p.Range "73.5%"
Bold-Metric $p.Range "`$4.7M"

# • Built real-time FEC analysis systems ... valued over $2 trillion
$p = $d.Paragraphs.Item(14)
Bold-Metric $p.Range "`$2"

# "Data Products Manager - Helm/Murmuration" bullet
# • Modernized legacy ETL processes ... reducing processing time by 57%
$p = $d.Paragraphs.Item(39)
Bold-Metric $p.Range "57%"

# "KEY ACHIEVEMENTS AND IMPACT" bullets
# • Algorithmic innovation: Pioneered trigonometric boundary estimation
#   reducing mapping costs 73.5%
$p = $d.Paragraphs.Item(55)
Bold-Metric $p.Range "73.5%"

# • $4.7M savings enabled nonprofit access
$p = $d.Paragraphs.Item(56)
Bold-Metric $p.Range "`$4.7M"

# • Platform impact: Built redistricting system serving 12,847 analysts
#   across 89 organizations
$p = $d.Paragraphs.Item(57)
Bold-Metric $p.Range "12,847"
